$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11 from 45184 to 45186
$ws.Range("C2:C11").Value2 = 45186

# Row 2 (A 21853-2022) - add friendly display text as second argument to HYPERLINK formulas
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/artfynd/A 21853-2022.xlsx", "A 21853-2022")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/kartor/A 21853-2022.png", "A 21853-2022")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/klagomål/A 21853-2022.docx", "A 21853-2022")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/klagomålsmail/A 21853-2022.docx", "A 21853-2022")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/tillsyn/A 21853-2022.docx", "A 21853-2022")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/tillsynsmail/A 21853-2022.docx", "A 21853-2022")'

# Row 3 (A 45510-2020) - add friendly display text as second argument to HYPERLINK formulas
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/artfynd/A 45510-2020.xlsx", "A 45510-2020")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/kartor/A 45510-2020.png", "A 45510-2020")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/klagomål/A 45510-2020.docx", "A 45510-2020")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/klagomålsmail/A 45510-2020.docx", "A 45510-2020")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/tillsyn/A 45510-2020.docx", "A 45510-2020")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_NACKA/tillsynsmail/A 45510-2020.docx", "A 45510-2020")'
